$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exclusions")

# Select the "ueiDUNS" column (column G) and remove it entirely; everything
# to its right shifts one column to the left.
$col = $ws.Columns("G")
$col.Select() | Out-Null
$col.Delete()
